$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values per the new TPM-based computation.
# Row 2
$ws.Range("G2").Value = 0.1806205
$ws.Range("H2").Value = 0.361241
$ws.Range("I2").Value = 0.7284509268949775
$ws.Range("J2").Value = 0.7284509268949775
$ws.Range("M2").Value = 31.6939195
$ws.Range("N2").Value = 63.387839
$ws.Range("O2").Value = 0.5590288178082639
$ws.Range("P2").Value = 0.5498517214736327
$ws.Range("Q2").Value = 5.724571587049749
$ws.Range("R2").Value = 22.898286348199
$ws.Range("S2").Value = 0.4072250604934333
$ws.Range("T2").Value = 0.4005399961622667

# Row 3
$ws.Range("G3").Value = 0.1806205
$ws.Range("H3").Value = 0.361241
$ws.Range("I3").Value = 0.7284509268949775
$ws.Range("J3").Value = 0.7284509268949775
$ws.Range("O3").Value = 0.01190832240406543
$ws.Range("P3").Value = 0.01756925054260132
$ws.Range("Q3").Value = 0.1219437029221666
$ws.Range("R3").Value = 0.7316622175329999
$ws.Range("S3").Value = 0.008674628493005692
$ws.Range("T3").Value = 0.01279833684260802

# Row 4
$ws.Range("G4").Value = 0.1806205
$ws.Range("H4").Value = 0.361241
$ws.Range("I4").Value = 0.7284509268949775
$ws.Range("J4").Value = 0.7284509268949775
$ws.Range("M4").Value = 0.6695323333333333
$ws.Range("N4").Value = 2.008597
$ws.Range("O4").Value = 0.01180945350693346
$ws.Range("P4").Value = 0.01742338176565341
$ws.Range("Q4").Value = 0.1209312648128333
$ws.Range("R4").Value = 0.7255875888769999
$ws.Range("S4").Value = 0.008602607353248821
$ws.Range("T4").Value = 0.01269207859683528

# Row 5
$ws.Range("G5").Value = 0.1806205
$ws.Range("H5").Value = 0.361241
$ws.Range("I5").Value = 0.7284509268949775
$ws.Range("J5").Value = 0.7284509268949775
$ws.Range("M5").Value = 23.1082075
$ws.Range("N5").Value = 46.216415
$ws.Range("O5").Value = 0.4075909235647884
$ws.Range("P5").Value = 0.4008998531735689
$ws.Range("Q5").Value = 4.17381599275375
$ws.Range("R5").Value = 16.695263971015
$ws.Range("S5").Value = 0.29690998606475
$ws.Range("T5").Value = 0.2920358696363466

# Row 6
$ws.Range("G6").Value = 0.1806205
$ws.Range("H6").Value = 0.361241
$ws.Range("I6").Value = 0.7284509268949775
$ws.Range("J6").Value = 0.7284509268949775
$ws.Range("M6").Value = 0.2922816666666667
$ws.Range("N6").Value = 0.876845
$ws.Range("O6").Value = 0.005155369773173549
$ws.Range("P6").Value = 0.007606107738040217
$ws.Range("Q6").Value = 0.05279206077416666
$ws.Range("R6").Value = 0.316752364645
$ws.Range("S6").Value = 0.003755433889754621
$ws.Range("T6").Value = 0.005540676231838457

# Row 7
$ws.Range("G7").Value = 0.1806205
$ws.Range("H7").Value = 0.361241
$ws.Range("I7").Value = 0.7284509268949775
$ws.Range("J7").Value = 0.7284509268949775
$ws.Range("M7").Value = 0.255529
$ws.Range("N7").Value = 0.7665869999999999
$ws.Range("O7").Value = 0.004507112942775279
$ws.Range("P7").Value = 0.006649685306503471
$ws.Range("Q7").Value = 0.04615377574449999
$ws.Range("R7").Value = 0.2769226544669999
$ws.Range("S7").Value = 0.003283210600785002
$ws.Range("T7").Value = 0.004843969425082366

# Row 8
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.067331
$ws.Range("H8").Value = 0.134662
$ws.Range("I8").Value = 0.2715490731050226
$ws.Range("J8").Value = 0.2715490731050226
$ws.Range("M8").Value = 31.6939195
$ws.Range("N8").Value = 63.387839
$ws.Range("O8").Value = 0.5590288178082639
$ws.Range("P8").Value = 0.5498517214736327
$ws.Range("Q8").Value = 2.1339832938545
$ws.Range("R8").Value = 8.535933175418
$ws.Range("S8").Value = 0.1518037573148306
$ws.Range("T8").Value = 0.149311725311366

# Row 9
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.067331
$ws.Range("H9").Value = 0.134662
$ws.Range("I9").Value = 0.2715490731050226
$ws.Range("J9").Value = 0.2715490731050226
$ws.Range("O9").Value = 0.01190832240406543
$ws.Range("P9").Value = 0.01756925054260132
$ws.Range("Q9").Value = 0.04545769423433333
$ws.Range("R9").Value = 0.272746165406
$ws.Range("S9").Value = 0.003233693911059743
$ws.Range("T9").Value = 0.004770913699993302

# Row 10
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.067331
$ws.Range("H10").Value = 0.134662
$ws.Range("I10").Value = 0.2715490731050226
$ws.Range("J10").Value = 0.2715490731050226
$ws.Range("M10").Value = 0.6695323333333333
$ws.Range("N10").Value = 2.008597
$ws.Range("O10").Value = 0.01180945350693346
$ws.Range("P10").Value = 0.01742338176565341
$ws.Range("Q10").Value = 0.04508028153566666
$ws.Range("R10").Value = 0.270481689214
$ws.Range("S10").Value = 0.003206846153684639
$ws.Range("T10").Value = 0.004731303168818137

# Row 11
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.067331
$ws.Range("H11").Value = 0.134662
$ws.Range("I11").Value = 0.2715490731050226
$ws.Range("J11").Value = 0.2715490731050226
$ws.Range("M11").Value = 23.1082075
$ws.Range("N11").Value = 46.216415
$ws.Range("O11").Value = 0.4075909235647884
$ws.Range("P11").Value = 0.4008998531735689
$ws.Range("Q11").Value = 1.5558987191825
$ws.Range("R11").Value = 6.22359487673
$ws.Range("S11").Value = 0.1106809375000384
$ws.Range("T11").Value = 0.1088639835372223

# Row 12
$ws.Range("F12").Value = 0.5
$ws.Range("G12").Value = 0.067331
$ws.Range("H12").Value = 0.134662
$ws.Range("I12").Value = 0.2715490731050226
$ws.Range("J12").Value = 0.2715490731050226
$ws.Range("M12").Value = 0.2922816666666667
$ws.Range("N12").Value = 0.876845
$ws.Range("O12").Value = 0.005155369773173549
$ws.Range("P12").Value = 0.007606107738040217
$ws.Range("Q12").Value = 0.01967961689833334
$ws.Range("R12").Value = 0.11807770139
$ws.Range("S12").Value = 0.001399935883418928
$ws.Range("T12").Value = 0.002065431506201761

# Row 13
$ws.Range("F13").Value = 0.5
$ws.Range("G13").Value = 0.067331
$ws.Range("H13").Value = 0.134662
$ws.Range("I13").Value = 0.2715490731050226
$ws.Range("J13").Value = 0.2715490731050226
$ws.Range("M13").Value = 0.255529
$ws.Range("N13").Value = 0.7665869999999999
$ws.Range("O13").Value = 0.004507112942775279
$ws.Range("P13").Value = 0.006649685306503471
$ws.Range("Q13").Value = 0.017205023099
$ws.Range("R13").Value = 0.103230138594
$ws.Range("S13").Value = 0.001223902341990278
$ws.Range("T13").Value = 0.001805715881421105
